# cs-en-us-009pct.xlsx weekly CompStat refresh
# - new Police Commissioner name
# - report header Volume/Number + week-covering dates bumped
# - precinct crime-stat figures refreshed for the new reporting week

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead -------------------------------------------------------------

# Police Commissioner name (M6)
$ws.Range("M6").Value = "Jessica S. Tisch"

# "Volume 31   Number  47" -> "... Number  48" (A8)
$a8 = $ws.Range("A8").Value2
$a8 = $a8.Replace("47", "48")
$ws.Range("A8").Value = $a8

# "Report Covering the Week  11/18/2024  Through  11/24/2024" (C9)
$c9 = $ws.Range("C9").Value2
$c9 = $c9.Replace("11/18/2024", "11/25/2024").Replace("11/24/2024", "12/1/2024")
$ws.Range("C9").Value = $c9

# --- Rape row (15): only the 2-Year % change moves ------------------------
$ws.Range("L15").Value = -5.555555555555

# --- Robbery row (16) -------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 17
$ws.Range("H16").Value = 41.666666666666
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 171
$ws.Range("K16").Value = -25.146198830409
$ws.Range("L16").Value = -38.755980861244
$ws.Range("M16").Value = -26.011560693641
$ws.Range("N16").Value = -85.964912280701

# --- Felony Assault row (17) -------------------------------------------------
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 17
$ws.Range("H17").Value = 13.333333333333
$ws.Range("I17").Value = 186
$ws.Range("J17").Value = 206
$ws.Range("K17").Value = -9.708737864077
$ws.Range("L17").Value = -3.125
$ws.Range("M17").Value = 31.914893617021
$ws.Range("N17").Value = -65.99634369287

# --- Burglary row (18) -------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 38
$ws.Range("H18").Value = -42.105263157894
$ws.Range("I18").Value = 156
$ws.Range("J18").Value = 262
$ws.Range("K18").Value = -40.458015267175
$ws.Range("L18").Value = -56.179775280898
$ws.Range("M18").Value = -24.271844660194
$ws.Range("N18").Value = -81.902552204176

# --- Grand Larceny row (19) --------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 11
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = -19.354838709677
$ws.Range("I19").Value = 700
$ws.Range("J19").Value = 910
$ws.Range("K19").Value = -23.076923076923
$ws.Range("L19").Value = -27.234927234927
$ws.Range("M19").Value = -1.269393511988
$ws.Range("N19").Value = -54.983922829582

# --- G.L.A. row (20) ----------------------------------------------------------
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 44
$ws.Range("J20").Value = 42
$ws.Range("K20").Value = 4.761904761904
$ws.Range("L20").Value = -2.222222222222
$ws.Range("M20").Value = -2.222222222222
$ws.Range("N20").Value = -91.489361702127

# --- TOTAL row (21) ------------------------------------------------------------
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -4.545454545454
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = -11.538461538461
$ws.Range("I21").Value = 1233
$ws.Range("J21").Value = 1610
$ws.Range("K21").Value = -23.416149068323
$ws.Range("L21").Value = -31.040268456375
$ws.Range("M21").Value = -4.418604651162
$ws.Range("N21").Value = -72.248480756245

# --- Transit row (22): WTD 2023/%chg go from "N/A" placeholders to real numbers
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = -8.333333333333
$ws.Range("M22").Value = -47.619047619047

# --- Housing row (23): same placeholder -> number transition ------------------
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = 100
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = 14.285714285714
$ws.Range("I23").Value = 84
$ws.Range("J23").Value = 116
$ws.Range("K23").Value = -27.586206896551
$ws.Range("L23").Value = -32.8
$ws.Range("M23").Value = -20.754716981132

# --- Petit Larceny row (24) -----------------------------------------------------
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -12.903225806451
$ws.Range("F24").Value = 144
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = 15.2
$ws.Range("I24").Value = 1449
$ws.Range("J24").Value = 1389
$ws.Range("K24").Value = 4.319654427645
$ws.Range("L24").Value = -28.655834564254
$ws.Range("M24").Value = -9.550561797752

# --- Retail Theft row (25) -------------------------------------------------------
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -35.294117647058
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 74
$ws.Range("H25").Value = 4.054054054054
$ws.Range("I25").Value = 833
$ws.Range("J25").Value = 785
$ws.Range("K25").Value = 6.114649681528
$ws.Range("L25").Value = -39.768618944323

# --- Misd. Assault row (26) -------------------------------------------------------
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 114.285714285714
$ws.Range("F26").Value = 45
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 9.756097560975
$ws.Range("I26").Value = 414
$ws.Range("J26").Value = 431
$ws.Range("K26").Value = -3.944315545243
$ws.Range("L26").Value = -3.720930232558
$ws.Range("M26").Value = -4.387990762124

# --- UCR Rape row (27): only 2-Year % change moves --------------------------------
$ws.Range("L27").Value = -39.393939393939

# --- Other Sex Crimes row (28) -----------------------------------------------------
# WTD-2024 (C28) goes from a real count back to the "N/A" placeholder -- copy
# the format+value straight from C22, which already holds that placeholder.
$ws.Range("C22").Copy($ws.Range("C28"))
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -62.5
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 1.666666666666
$ws.Range("L28").Value = -29.885057471264

# --- Shooting Vic. row (29): WTD-2024 (C29) -> "N/A" placeholder ------------------
$ws.Range("C22").Copy($ws.Range("C29"))

# --- Shooting Inc. row (30): WTD-2024 (C30) -> "N/A" placeholder ------------------
$ws.Range("C22").Copy($ws.Range("C30"))

# --- Hate Crimes row (33): 28-Day % change column -----------------------------------
$ws.Range("L33").Value = 0
